{"js": "// Remove the \"Author\" byline paragraphs (\"Ben Jarman\" and \"Catherine Heard\")\n// from the document, leaving the Title/Subtitle/Date/Abstract etc. untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style,text\");\nawait context.sync();\n\n// Collect the paragraphs to remove first, then delete from last to first so\n// that earlier deletions don't shift the indices/ranges of later ones.\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Author\") {\n    toDelete.push(para);\n  }\n}\n\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Author\" byline paragraphs (\"Ben Jarman\" and \"Catherine Heard\")\n# from the document, leaving the Title/Subtitle/Date/Abstract etc. untouched.\n$doc = $word.ActiveDocument\n\n# Snapshot the paragraphs collection into an array first, then walk it\n# back-to-front so deleting one paragraph doesn't invalidate/shift the\n# ranges of the ones still to be processed.\n$all = @($doc.Paragraphs)\nfor ($i = $all.Count - 1; $i -ge 0; $i--) {\n    $p = $all[$i]\n    if ($p.Style.NameLocal -eq \"Author\") {\n        $p.Range.Delete()\n    }\n}\n"}
